# Edit script: apply "VSC Trunk Name" + "Segmentation Id" fields to the Vscs sheet
# of the OpenStack install workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# --- 1. Insert the new "VSC Trunk Name" row above the old row 54 -----------
$ws.Range("A54").EntireRow.Insert()

# --- 2. Insert the three new "Segmentation Id" rows (now at 75-77) --------
$ws.Range("A75:A77").EntireRow.Insert()

# --- 3. Fix up cell formatting for the newly inserted rows -----------------
# Column A keeps the "label" style copied down from the row above (style 6).
# Columns B/C need the "input" style (style 7) instead of the copied style 6.
$ws.Range("B55:C55").Copy()
$ws.Range("B54:C54").PasteSpecial(-4122)

$ws.Range("B74:C74").Copy()
$ws.Range("B75:C77").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 4. Set the label text for the newly inserted rows ----------------------
$ws.Range("A54").Value = "VSC Trunk Name"
$ws.Range("A75").Value = "First External Port's Segmentation Id"
$ws.Range("A76").Value = "Second External Port's Segmentation Id"
$ws.Range("A77").Value = "Third External Port's Segmentation Id"

# --- 5. Rebuild the comments ------------------------------------------------
# Row insertion does not relocate existing cell comments, so every comment
# on this sheet is cleared and re-added at its correct final location.
$ws.Cells.ClearComments()

$ws.Range("A5").AddComment("Hostname of the VSC instance")
$ws.Range("A6").AddComment("The BOF/Management IP address of the controller")
$ws.Range("A7").AddComment("Management network prefix length")
$ws.Range("A8").AddComment("Gateway IP on the Management network")
$ws.Range("A9").AddComment("Static Route list to be configured on the management/BOF interface. Define as empty list if no static routes are to be configured. [default: [ 0.0.0.0/1, 128.0.0.0/1 ]] (List items separated by comma.)")
$ws.Range("A11").AddComment("The Control/Data IP address of the controller")
$ws.Range("A12").AddComment("Control network prefix length")
$ws.Range("A13").AddComment("The VLAN ID for the uplink connection of the controller. This field is deprecated and will be removed in future releases. Use ctrl_ip_vprn_list instead. [default: 0]")
$ws.Range("A14").AddComment("The Control/Data VPRN IP address of the controller. This field is deprecated and will be removed in future releases. Use ctrl_ip_vprn_list instead.")
$ws.Range("A15").AddComment("Comma Seperated list of control IP VPRNs. Each item in list is of format VLAN_ID:Control_IP_ADDR/NETMASK_LENGTH:NEXT_HOP . Control IP address, netmask length and next hop are optional (List items separated by comma.)")
$ws.Range("A17").AddComment("Required for BGP pairing with peers [default: ]")
$ws.Range("A18").AddComment("The router ID of this VSC in IPV4 address format. Required when system_ip is IPV6. [default: (System IP)]")
$ws.Range("A19").AddComment("Name of the virtual machine on the Hypervisor/vCenter Server. [default: (Hostname)]")
$ws.Range("A20").AddComment("Name of the system if other than hostname [default: (Hostname)]")
$ws.Range("A21").AddComment("Unique username used to identify this VSC in its XMPP connection with VSD [default: vsc1]")
$ws.Range("A23").AddComment("Type of hypervisor environment where VMs will be instantiated. Use 'none' when skipping predeploy.")
$ws.Range("A24").AddComment("Hostname or IP address of the hypervisor where VM  will be instantiated. In the case of deployment in a vCenter environment, this will be the FQDN of the vCenter Server. When target_server_type is set to openstack, this property is unused and should be set to 0.0.0.0")
$ws.Range("A26").AddComment("Network Bridge used for the management interface of a component or the BOF interface on VSC. This will be a Distributed Virtual PortGroup (DVPG) when deploying on vCenter or a Linux network bridge when deploying on KVM. This field can be overridden by defining the management network bridge separately in the component configuration. Defaults to the global setting [default: (global Bridge interface)]")
$ws.Range("A27").AddComment("Network Bridge used for the data path of a component or the Control interface on VSC. This will be a Distributed Virtual PortGroup (DVPG) when deploying on vCenter or a Linux network bridge when deploying on KVM. [default: (global Bridge interface)]")
$ws.Range("A28").AddComment("FQDN of the VSD or VSD cluster for this VSC")
$ws.Range("A29").AddComment("Private Management IP Address of VSC instances")
$ws.Range("A30").AddComment("Private Control IP Address of VSC Instances")
$ws.Range("A31").AddComment("Private Data Gateway IP Address of VSC Instances")
$ws.Range("A32").AddComment("List of route reflector IP addresses if present (List items separated by comma.)")
$ws.Range("A34").AddComment("IP Address for Optional BGP Interface")
$ws.Range("A35").AddComment("Prefix length for the optional BGP interface [default: 24]")
$ws.Range("A36").AddComment("VLAN ID for the optional BGP interface [default: 1000]")
$ws.Range("A38").AddComment("Name of the vCenter Datacenter on which the VSC VM will be deployed. Defaults to the common vCenter Datacenter Name if not defined here. [default: (global vCenter Datacenter Name)]")
$ws.Range("A39").AddComment("Name of the vCenter Cluster on which the VSC VM will be deployed. Defaults to the common vCenter Cluster Name if not defined here. [default: (global vCenter Cluster Name)]")
$ws.Range("A40").AddComment("Requires ovftool 4.3. Reference to the host on the vCenter cluster on which to deploy Nuage components [default: (global vCenter Host Reference)]")
$ws.Range("A41").AddComment("Name of the vCenter Datastore on which the VSC VM will be deployed. Defaults to the common vCenter Datastore Name if not defined here. [default: (global vCenter Datastore Name)]")
$ws.Range("A42").AddComment("Optional path to a folder defined on vCenter where VM will be instantiated [default: (global vCenter VM folder)]")
$ws.Range("A43").AddComment("Optional path to a hosts and clusters folder defined on vCenter where VM will be instantiated")
$ws.Range("A45").AddComment("Name of image installed on OpenStack for VSC")
$ws.Range("A46").AddComment("Name of instance flavor installed on OpenStack for VSC")
$ws.Range("A47").AddComment("Name of availability zone on OpenStack for VSC")
$ws.Range("A48").AddComment("Name of management network on OpenStack for VSC")
$ws.Range("A49").AddComment("Name of management subnet on OpenStack for VSC")
$ws.Range("A50").AddComment("Name for Mgmt interface")
$ws.Range("A51").AddComment("Set of security groups to associate with Mgmt interface (List items separated by comma.)")
$ws.Range("A52").AddComment("Name of control network on OpenStack for VSC")
$ws.Range("A53").AddComment("Name of control subnet on OpenStack for VSC")
$ws.Range("A54").AddComment("Name of the trunk to be setup between control port and its underlay ports")
$ws.Range("A55").AddComment("Name of first external network on OpenStack for VSC")
$ws.Range("A56").AddComment("Name of first external subnet on OpenStack for VSC")
$ws.Range("A57").AddComment("Name of second external network on OpenStack for VSC")
$ws.Range("A58").AddComment("Name of second external subnet on OpenStack for VSC")
$ws.Range("A59").AddComment("Name of third external network on OpenStack for VSC")
$ws.Range("A60").AddComment("Name of third external subnet on OpenStack for VSC")
$ws.Range("A61").AddComment("Name for Control interface")
$ws.Range("A62").AddComment("Set of security groups to associate with Control interface (List items separated by comma.)")
$ws.Range("A63").AddComment("Name for first external port interface")
$ws.Range("A64").AddComment("Set of security groups to associate with first external interface (List items separated by comma.)")
$ws.Range("A65").AddComment("Name for second external port interface")
$ws.Range("A66").AddComment("Set of security groups to associate with second external interface (List items separated by comma.)")
$ws.Range("A67").AddComment("Name for third external port interface")
$ws.Range("A68").AddComment("Set of security groups to associate with first external interface (List items separated by comma.)")
$ws.Range("A69").AddComment("The first External IP address of the controller")
$ws.Range("A70").AddComment("First External network prefix length")
$ws.Range("A71").AddComment("The second External IP address of the controller")
$ws.Range("A72").AddComment("Second External network prefix length")
$ws.Range("A73").AddComment("The third External IP address of the controller")
$ws.Range("A74").AddComment("Third External network prefix length")
$ws.Range("A75").AddComment("Segmentation id of first external port to be used during trunking")
$ws.Range("A76").AddComment("Segmentation id of first external port to be used during trunking")
$ws.Range("A77").AddComment("Segmentation id of third external port to be used during trunking")
$ws.Range("A78").AddComment("Name for Mgmt interface")
$ws.Range("A79").AddComment("Set of security groups to associate with Mgmt interface (List items separated by comma.)")
$ws.Range("A81").AddComment("Used in postdeploy and health workflows as expected values if non-zero [default: 0]")
$ws.Range("A82").AddComment("Used in postdeploy and health workflows as expected values if non-zero [default: 0]")
$ws.Range("A83").AddComment("Used in postdeploy and health workflows as expected values if non-zero [default: 0]")
$ws.Range("A84").AddComment("Used in postdeploy and health workflows as expected values if non-zero [default: 0]")
$ws.Range("A85").AddComment("Used in postdeploy and health workflows as expected values if non-zero [default: 0]")
$ws.Range("A87").AddComment("Ejabberd user id used to create the certificate")
$ws.Range("A88").AddComment("Path to VSC certificate key pem file")
$ws.Range("A89").AddComment("Path to VSC certificate pem file")
$ws.Range("A90").AddComment("Path to CA certificate pem file")
$ws.Range("A91").AddComment("XMPP domain used in custom certificates")
$ws.Range("A92").AddComment("Name of the credentials set for the vsc")
$ws.Range("A94").AddComment("Cpuset information for cpu pinning on KVM. For example, VSC requires 4 cores and sample values will be of the form [ 0, 1, 2, 3 ] (List items separated by comma.)")
$ws.Range("A95").AddComment("Enables hardening configuration on VSC [default: True]")
$ws.Range("A96").AddComment("Paths to files that can be optionally applied for additional VSC configuration (List items separated by comma.)")
$ws.Range("A97").AddComment("This will override the Metro Provided config on the VSC by the config provided in vsc_config_file_paths [default: False]")

Write-Output "Vscs sheet updated: VSC Trunk Name + Segmentation Id rows inserted."
